$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1. Paragraphs that only need <w:ilvl w:val="0"/> inserted into their
#    existing <w:numPr> (numId stays "0"); text/runs are unchanged.
# ---------------------------------------------------------------------------

# Paragraph 5: "git branch test"
$xml5 = @"
<w:p $ns><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>git branch test</w:t></w:r></w:p>
"@
$d.Paragraphs(5).Range.InsertXML($xml5)

# Paragraph 7: "git checkout test"
$xml7 = @"
<w:p $ns><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>git checkout test</w:t></w:r></w:p>
"@
$d.Paragraphs(7).Range.InsertXML($xml7)

# Paragraph 9: "git merge test"
$xml9 = @"
<w:p $ns><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>git merge test</w:t></w:r></w:p>
"@
$d.Paragraphs(9).Range.InsertXML($xml9)

# Paragraph 11: "git breach -d test"
$xml11 = @"
<w:p $ns><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>git breach -d test</w:t></w:r></w:p>
"@
$d.Paragraphs(11).Range.InsertXML($xml11)

# ---------------------------------------------------------------------------
# 2. Paragraph 13 ("git push origin 本地分支名"): add <w:ilvl>, flip the
#    paragraph-mark rFonts hint from "default" to "eastAsia", drop the
#    _GoBack bookmark from here (it moves further down), and append the new
#    paragraphs describing the push-failure / pull-then-push workaround plus
#    two trailing blank paragraphs -- all in one InsertXML call so the
#    paragraph count grows by exactly nine.
# ---------------------------------------------------------------------------

$xml13 = @"
<w:p $ns><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>git push origin 本地分支名</w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="0" w:leftChars="0" w:firstLine="0" w:firstLineChars="0"/><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>本地仓库push到远程仓库失败，解决办法</w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>先pull使本地仓库与远程仓库依赖一致，然后再push到远程仓库，步骤：</w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>先pull:</w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve">git pull origin master  </w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>后push:</w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>git push -u origin master</w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
<w:p><w:pPr><w:numPr><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:numPr><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cstheme="minorEastAsia"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr></w:p>
"@
$d.Paragraphs(13).Range.InsertXML($xml13)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
